$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.104696273803711
$ws.Range("B1").Value = 2.49763011932373
$ws.Range("C1").Value = 1.847217440605164
$ws.Range("D1").Value = 1.72240686416626
$ws.Range("E1").Value = 1.644418597221375
